# Update the "Estado de Cuenta" worker table (rows 16-20) on Hoja1.
# Previous EC records are removed and replaced with new ones (reordered),
# and the underlying "database" (Valor Mora / Salario Basico) values are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Target state for columns C (N Doc Trabajador), D (Nombre Trabajador),
# F (Valor Mora) and G (Salario Basico) for rows 16-20.
$data = @(
    @{ Row = 16; Id = "1143361023"; Name = "MIGUEL REINIER GARCIA SANJUAN"; F = 200000; G = 5000000 },
    @{ Row = 17; Id = "1143360831"; Name = "LIS LAURY CASTILLO MARTINEZ";   F = 80000;  G = 2000000 },
    @{ Row = 18; Id = "1047387594"; Name = "DANILO YESITH BARRIOS CANAVAL"; F = 64000;  G = 1600000 },
    @{ Row = 19; Id = "1143357850"; Name = "VANESSA ACEVEDO VILLEROS";      F = 64000;  G = 1423500 },
    @{ Row = 20; Id = "45527359";   Name = "ANA OLINDA HERNANDEZ PINTO";    F = 46400;  G = 1160000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Id
    $ws.Cells.Item($r, 4).Value = $item.Name
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
